$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column is treated as text so values like "0.100" or
# "7.90" keep their exact (dotted-thousands / trailing-zero) formatting
# instead of being auto-coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "35.078.02"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "1.809.01"
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "232.82"
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "40.47"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("E9").Value = "  +5.81%  "
$ws.Range("D10").Value = "0.0684"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").Value = "0.100"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "2.070.54"
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("D13").Value = "1.811.42"
$ws.Range("E13").Value = "  -2.01%  "
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").Value = "11.05"
$ws.Range("E15").Value = "  -5.03%  "
$ws.Range("D16").Value = "4.66"
$ws.Range("E16").Value = "  -1.51%  "
$ws.Range("D17").Value = "35.024.97"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").Value = "69.71"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").Value = "237.66"
$ws.Range("E20").Value = "  -3.24%  "
$ws.Range("D21").Value = "11.95"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").Value = "2.25"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("D25").Value = "172.01"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("D26").Value = "7.90"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "17.51"
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("E28").Value = "  -1.70%  "
$ws.Range("E29").Value = "  +17.74%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  +4.37%  "
$ws.Range("D32").Value = "0.0558"
$ws.Range("E32").Value = "  +4.17%  "
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("D34").Value = "1.77"
$ws.Range("E34").Value = "  -6.22%  "
$ws.Range("E35").Value = "  +3.35%  "
$ws.Range("E36").Value = "  +5.77%  "
$ws.Range("D37").Value = "92.35"
$ws.Range("E37").Value = "  +3.82%  "
$ws.Range("D38").Value = "0.0194"
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("D39").Value = "1.315.91"
$ws.Range("E39").Value = "  -1.91%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -3.80%  "
$ws.Range("D42").Value = "2.46"
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("D43").Value = "14.45"
$ws.Range("E43").Value = "  -2.95%  "
$ws.Range("D44").Value = "2.25"
$ws.Range("E44").Value = "  -6.95%  "
$ws.Range("E45").Value = "  -2.14%  "
$ws.Range("D46").Value = "6.31"
$ws.Range("E46").Value = "  +4.23%  "
$ws.Range("D47").Value = "0.0512"
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("D48").Value = "1.987.42"
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").Value = "0.0667"
$ws.Range("E50").Value = "  +6.89%  "
$ws.Range("D51").Value = "99.46"
$ws.Range("E51").Value = "  -4.87%  "
